$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 668033.75
$ws.Range("I28").Value = 1122.625
$ws.Range("J28").Value = 1430217.9
$ws.Range("K28").Value = 1122.625
$ws.Range("L28").Value = 1430217.9
$ws.Range("M28").Value = -637.625
$ws.Range("N28").Value = -1431187.9

$ws.Range("H64").Value = 4374.0625
$ws.Range("I64").Value = 4416.706
$ws.Range("J64").Value = 4325.7334
$ws.Range("K64").Value = 4416.706
$ws.Range("L64").Value = 4325.7334
$ws.Range("M64").Value = -4168.706
$ws.Range("N64").Value = -4821.7334

$ws.Range("H67").Value = 4374.0625
$ws.Range("I67").Value = 4416.706
$ws.Range("J67").Value = 4325.7334
$ws.Range("K67").Value = 4416.706
$ws.Range("L67").Value = 4325.7334
$ws.Range("M67").Value = -3558.706
$ws.Range("N67").Value = -6041.7334

$ws.Range("H98").Value = 2453.4546
$ws.Range("I98").Value = 1698.8
$ws.Range("K98").Value = 1698.8
$ws.Range("M98").Value = -200.8

$ws.Range("H112").Value = 3810.2222
$ws.Range("I112").Value = 2100
$ws.Range("K112").Value = 6300
$ws.Range("M112").Value = -5192

$ws.Range("H116").Value = 4066.2307
$ws.Range("I116").Value = 3079.5
$ws.Range("K116").Value = 3079.5
$ws.Range("M116").Value = 362.5

$ws.Range("H122").Value = 2453.4546
$ws.Range("I122").Value = 1698.8
$ws.Range("K122").Value = 5096.4
$ws.Range("M122").Value = -2646.4

$ws.Range("H132").Value = 2574.7258
$ws.Range("I132").Value = 2323.5615
$ws.Range("K132").Value = 6970.684499999999
$ws.Range("M132").Value = -4440.684499999999

$ws.Range("H137").Value = 6919.216
$ws.Range("I137").Value = 5840.0713
$ws.Range("K137").Value = 17520.2139
$ws.Range("M137").Value = -14970.2139

$ws.Range("H138").Value = 5480.4443
$ws.Range("I138").Value = 3926.923
$ws.Range("J138").Value = 5884.36
$ws.Range("K138").Value = 11780.769
$ws.Range("L138").Value = 17653.08
$ws.Range("M138").Value = -6640.769
$ws.Range("N138").Value = -27933.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 154.8
$ws.Range("I4").Value = 131
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 131
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = -15
$ws.Range("N4").Value = -482

$ws.Range("H32").Value = 6526.206
$ws.Range("I32").Value = 5512.0303
$ws.Range("K32").Value = 5512.0303
$ws.Range("M32").Value = -5225.0303

$ws.Range("H61").Value = 20003970
$ws.Range("I61").Value = 2874
$ws.Range("J61").Value = 166678670
$ws.Range("K61").Value = 2874
$ws.Range("L61").Value = 166678670
$ws.Range("M61").Value = -2662
$ws.Range("N61").Value = -166679094

$ws.Range("H97").Value = 1927.6154
$ws.Range("I97").Value = 2198.4546
$ws.Range("K97").Value = 2198.4546
$ws.Range("M97").Value = -1702.4546

$ws.Range("H136").Value = 20003970
$ws.Range("I136").Value = 2874
$ws.Range("J136").Value = 166678670
$ws.Range("K136").Value = 8622
$ws.Range("L136").Value = 500036010
$ws.Range("M136").Value = -6072
$ws.Range("N136").Value = -500041110

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5026.3184
$ws.Range("I86").Value = 2992.5625
$ws.Range("J86").Value = 10449.667
$ws.Range("K86").Value = 2992.5625
$ws.Range("L86").Value = 10449.667
$ws.Range("M86").Value = -1869.5625
$ws.Range("N86").Value = -12695.667

$ws.Range("H89").Value = 5026.3184
$ws.Range("I89").Value = 2992.5625
$ws.Range("J89").Value = 10449.667
$ws.Range("K89").Value = 14962.8125
$ws.Range("L89").Value = 52248.335
$ws.Range("M89").Value = -9346.8125
$ws.Range("N89").Value = -63480.335

$ws.Range("H94").Value = 1629.1666
$ws.Range("I94").Value = 913.5
$ws.Range("J94").Value = 3060.5
$ws.Range("K94").Value = 913.5
$ws.Range("L94").Value = 3060.5
$ws.Range("M94").Value = -462.5
$ws.Range("N94").Value = -3962.5

$ws.Range("H134").Value = 2067.92
$ws.Range("I134").Value = 1457.5
$ws.Range("J134").Value = 3637.5715
$ws.Range("K134").Value = 4372.5
$ws.Range("L134").Value = 10912.7145
$ws.Range("M134").Value = -1837.5
$ws.Range("N134").Value = -15982.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 83340000
$ws.Range("I62").Value = 9997.5
$ws.Range("J62").Value = 250000000
$ws.Range("K62").Value = 9997.5
$ws.Range("L62").Value = 250000000
$ws.Range("M62").Value = -9373.5
$ws.Range("N62").Value = -250001248

$ws.Range("H65").Value = 83340000
$ws.Range("I65").Value = 9997.5
$ws.Range("J65").Value = 250000000
$ws.Range("K65").Value = 49987.5
$ws.Range("L65").Value = 1250000000
$ws.Range("M65").Value = -46867.5
$ws.Range("N65").Value = -1250006240

$ws.Range("H107").Value = 503.6129
$ws.Range("I107").Value = 328.1905
$ws.Range("K107").Value = 328.1905
$ws.Range("M107").Value = 1591.8095

$ws.Range("H132").Value = 4015.359
$ws.Range("I132").Value = 3210.3845
$ws.Range("K132").Value = 9631.1535
$ws.Range("M132").Value = -7101.1535

$ws.Range("H134").Value = 2759.9756
$ws.Range("I134").Value = 1841.5834
$ws.Range("J134").Value = 9372.4
$ws.Range("K134").Value = 5524.7502
$ws.Range("L134").Value = 28117.2
$ws.Range("M134").Value = -2989.7502
$ws.Range("N134").Value = -33187.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1942.8334
$ws.Range("I5").Value = 1942.8334
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5828.5002
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -5716.5002

$ws.Range("H68").Value = 2930.3262
$ws.Range("J68").Value = 3049.1428
$ws.Range("L68").Value = 9147.428400000001
$ws.Range("N68").Value = -10769.4284

$ws.Range("H71").Value = 2930.3262
$ws.Range("J71").Value = 3049.1428
$ws.Range("L71").Value = 27442.2852
$ws.Range("N71").Value = -35554.2852

$ws.Range("H81").Value = 1994.25
$ws.Range("I81").Value = 1994.25
$ws.Range("K81").Value = 5982.75
$ws.Range("M81").Value = -4859.75

$ws.Range("H84").Value = 1994.25
$ws.Range("I84").Value = 1994.25
$ws.Range("K84").Value = 17948.25
$ws.Range("M84").Value = -12332.25

$ws.Range("H107").Value = 1647.2
$ws.Range("I107").Value = 1542.2222
$ws.Range("J107").Value = 1733.091
$ws.Range("K107").Value = 4626.6666
$ws.Range("L107").Value = 5199.272999999999
$ws.Range("M107").Value = -2706.6666
$ws.Range("N107").Value = -9039.272999999999

$ws.Range("H113").Value = 885.5714
$ws.Range("J113").Value = 956.4167
$ws.Range("L113").Value = 2869.2501
$ws.Range("N113").Value = -7209.2501

$ws.Range("H132").Value = 3384.6667
$ws.Range("I132").Value = 3762.6
$ws.Range("J132").Value = 3114.7144
$ws.Range("K132").Value = 33863.4
$ws.Range("L132").Value = 28032.4296
$ws.Range("M132").Value = -31333.4
$ws.Range("N132").Value = -33092.4296

$ws.Range("H135").Value = 1942.8334
$ws.Range("I135").Value = 1942.8334
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 17485.5006
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -14950.5006

$ws.Range("H136").Value = 3151.6155
$ws.Range("I136").Value = 2485.889
$ws.Range("J136").Value = 4649.5
$ws.Range("K136").Value = 7457.667
$ws.Range("L136").Value = 13948.5
$ws.Range("M136").Value = -2357.667
$ws.Range("N136").Value = -24148.5

$ws.Range("H140").Value = 1595.1177
$ws.Range("I140").Value = 1341.75
$ws.Range("J140").Value = 2203.2
$ws.Range("K140").Value = 4025.25
$ws.Range("L140").Value = 6609.599999999999
$ws.Range("M140").Value = 1154.75
$ws.Range("N140").Value = -16969.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 502
$ws.Range("I97").Value = 502
$ws.Range("K97").Value = 502
$ws.Range("M97").Value = -6

$ws.Range("H122").Value = 2594.889
$ws.Range("I122").Value = 2343.7188
$ws.Range("J122").Value = 3213.1538
$ws.Range("K122").Value = 7031.1564
$ws.Range("L122").Value = 9639.4614
$ws.Range("M122").Value = -4581.1564
$ws.Range("N122").Value = -14539.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 780.129
$ws.Range("I55").Value = 380.4375
$ws.Range("K55").Value = 380.4375
$ws.Range("M55").Value = -207.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 11091.192
$ws.Range("I96").Value = 1362.4286
$ws.Range("K96").Value = 1362.4286
$ws.Range("M96").Value = 10.57140000000004

$ws.Range("H107").Value = 690710.25
$ws.Range("I107").Value = 834240.7
$ws.Range("J107").Value = 1764.2
$ws.Range("K107").Value = 2502722.1
$ws.Range("L107").Value = 5292.6
$ws.Range("M107").Value = -2500802.1
$ws.Range("N107").Value = -9132.6
